# The commit replaces the "Managed daily production..." bullet with a new
# DocuSign/Assetmark bullet, and adds a brand-new "Developed and
# orchestrated Function Apps..." bullet right after the "Optimized site
# performance..." bullet (i.e. at the end of the first "Incedo" job entry).

$d = $word.ActiveDocument

# 1. Replace the "Managed daily production..." bullet text in place.
$d.Content.Find.Execute(
    "Managed daily production and non-production deployments, resolving non-production issues promptly to ensure smooth operations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implemented and Integrated DocuSign feature for Assetmark Tracking Center to streamline electronic signature process",
    2) | Out-Null

# 2. Insert a brand-new list-bullet paragraph right after the "Optimized
#    site performance..." bullet, matching its list/paragraph formatting
#    (Word does this automatically when splitting a paragraph).
$rng = $d.Content
$rng.Find.Execute(
    "boosted user satisfaction by 60%",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $rng.Paragraphs(1)
$anchorStart = $anchorPara.Range.Start
$anchorEnd = $anchorPara.Range.End

$anchorPara.Range.InsertParagraphAfter()

$allParas = $d.Paragraphs
for ($i = 1; $i -le $allParas.Count; $i++) {
    $p = $allParas.Item($i)
    if ($p.Range.Start -eq $anchorStart -and $p.Range.End -eq $anchorEnd) {
        $newPara = $allParas.Item($i + 1)
        $newPara.Range.Text = "Developed and orchestrated Function Apps and Logic Apps for APIs on Azure"
        break
    }
}
